$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 82, shifting existing rows 82-84 down to 83-85.
$ws.Rows.Item(82).Insert()

# Populate the newly inserted row 82 with the new weekly data entry.
$ws.Cells.Item(82, 1).Value = 1
$ws.Cells.Item(82, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(82, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(82, 4).Value = 45267
$ws.Cells.Item(82, 4).NumberFormat = $ws.Cells.Item(83, 4).NumberFormat
$ws.Cells.Item(82, 5).Value = 15
$ws.Cells.Item(82, 6).Value = 100112052
$ws.Cells.Item(82, 7).Value = "Albahaca"
$ws.Cells.Item(82, 8).Value = "Sin especificar"
$ws.Cells.Item(82, 9).Value = "Primera"
$ws.Cells.Item(82, 10).Value = 300
$ws.Cells.Item(82, 11).Value = 1000
$ws.Cells.Item(82, 12).Value = 1200
$ws.Cells.Item(82, 13).Value = 1100
$ws.Cells.Item(82, 14).Value = "$/paquete"
$ws.Cells.Item(82, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(82, 16).Value = 1100
$ws.Cells.Item(82, 17).Value = 1
$ws.Cells.Item(82, 18).Value = "Hortaliza"
